# TAC-3831 Fix enable TMS to import trips and edit lists of sub category in excel file
# Adds two new columns (Pickup Facility / Drop off Facility) to the shipment import template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in columns H and I
$ws.Cells.Item(1, 8).Value = "Pickup Facility  منطقة الانطلاق"
$ws.Cells.Item(1, 9).Value = "Drop off Facility   منطقة الوصول"

# Match the column widths used for the new columns
$ws.Columns.Item(8).ColumnWidth = 28.63
$ws.Columns.Item(9).ColumnWidth = 26.0

# Move the active selection to A2
$ws.Range("A2").Select()
